$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.692.62"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "3.328.73"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.08%  "
$ws.Range("D9").Value = "3.323.81"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.181"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.581"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.49"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "693.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "3.868.18"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.44"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.09%  "
$ws.Range("D17").Value = "67.648.58"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("E18").Value = "  -0.81%  "
$ws.Range("D19").Value = "3.316.75"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.894"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.98"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.48"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "567.87"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("E36").Value = "  -0.12%  "
$ws.Range("D37").Value = "3.712.00"
$ws.Range("E37").Value = "  -4.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.131"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.20"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").Value = "0.0₃0674"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0405"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.129"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.70%  "
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("E50").Value = "  -4.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.15%  "
